$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new review row (row 11) by copying row 10's formatting/types,
# then overwrite the cells that actually change content.
$ws.Range("A10:G10").Copy($ws.Range("A11:G11"))

# New comment text for the existing last row (row 10), column F was empty.
$ws.Range("F10").Value = "Great app Great story"

# Row 11 new review data.
$ws.Range("C11").Value = "cohenyossi408@gmail.com"
$ws.Range("D11").Value = "cohenn167@gmail.com"
$ws.Range("F11").Value = "Everybody need to try this game. Very fun and entertainning"

# Leave the selection on the newly edited cells, as in the authored file.
$ws.Range("C11:D11").Select()
